$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows (21..73) down to (22..74)
$ws.Range("A21").EntireRow.Insert()

# Populate the newly inserted row 21 with the new weekly record
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 44953
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100103
$ws.Range("H21").Value = "Frutos de hueso (carozo)"
$ws.Range("I21").Value = 100103002
$ws.Range("J21").Value = "Ciruela"
$ws.Range("K21").Value = "Black Amber"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 11000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 11500
$ws.Range("Q21").Value = "$/bandeja 18 kilos granel"
$ws.Range("R21").Value = "Provincia de Curicó"
$ws.Range("S21").Value = 639
$ws.Range("T21").Value = 18
